# Revert "AIP-400 AIP-443 Removed TestData File"
# This reverts commit 73d3c63c003632d20007a9ce5b0448c85efa5238, restoring the
# original BTC_95_BTC106.xlsx test-data sheet: the "Admin" credential row
# (row 7) that had been blanked out, and the nine timezone display names
# that had been swapped for non-standard variants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 lost its Username/Password ("Admin") values - restore them, matching
# every other data row (2-76).
$ws.Cells.Item(7, 1).Value = "Admin"
$ws.Cells.Item(7, 2).Value = "Admin"

# Restore the original timezone display names (column G) that were replaced
# with alternate/abbreviated spellings.
$ws.Cells.Item(2, 7).Value  = "India Standard Time"
$ws.Cells.Item(8, 7).Value  = "Mexico Standard Time 2"
$ws.Cells.Item(9, 7).Value  = "U.S. Mountain Standard Time"
$ws.Cells.Item(12, 7).Value = "Mexico Standard Time"
$ws.Cells.Item(15, 7).Value = "U.S. Eastern Standard Time"
$ws.Cells.Item(16, 7).Value = "S.A. Pacific Standard Time"
$ws.Cells.Item(18, 7).Value = "S.A. Western Standard Time"
$ws.Cells.Item(19, 7).Value = "Pacific S.A. Standard Time"
$ws.Cells.Item(20, 7).Value = "Newfoundland and Labrador Standard Time"

# Restore the original selection (cell G2) shown in the workbook view.
$ws.Range("G2").Select() | Out-Null
